$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.956.66'
$ws.Range("E2").Value = '  -3.47%  '

$ws.Range("D3").Value = '3.316.55'
$ws.Range("E3").Value = '  -5.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '183.00'
$ws.Range("E5").Value = '  -7.03%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '532.71'
$ws.Range("E6").Value = '  -2.90%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.606'
$ws.Range("E7").Value = '  +0.76%  '

$ws.Range("D8").Value = '3.314.37'
$ws.Range("E8").Value = '  -4.85%  '

$ws.Range("E9").Value = '  +0.07%  '

$ws.Range("E10").Value = '  -4.51%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '60.15'
$ws.Range("E11").Value = '  -2.95%  '

$ws.Range("E12").Value = '  -5.46%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000263'
$ws.Range("E13").Value = '  -1.53%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.17'
$ws.Range("E14").Value = '  -5.87%  '

$ws.Range("D15").Value = '3.858.04'
$ws.Range("E15").Value = '  -4.73%  '

$ws.Range("D16").Value = '3.323.94'
$ws.Range("E16").Value = '  -4.68%  '

$ws.Range("E17").Value = '  -4.69%  '

$ws.Range("D19").Value = '64.786.69'
$ws.Range("E19").Value = '  -3.26%  '

$ws.Range("E20").Value = '  -4.52%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.965'
$ws.Range("E21").Value = '  -5.44%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '377.69'
$ws.Range("E22").Value = '  -2.49%  '

$ws.Range("E23").Value = '  -3.66%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.68'
$ws.Range("E24").Value = '  -0.36%  '

$ws.Range("E25").Value = '  -4.14%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.88'
$ws.Range("E26").Value = '  +3.94%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.08'
$ws.Range("E27").Value = '  -0.97%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.71'
$ws.Range("E28").Value = '  -2.46%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '11.69'
$ws.Range("E29").Value = '  -3.28%  '

$ws.Range("E30").Value = '  -2.23%  '

$ws.Range("E31").Value = '  -5.28%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '657.56'
$ws.Range("E32").Value = '  -2.67%  '

$ws.Range("E33").Value = '  -1.96%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.38'
$ws.Range("E34").Value = '  -2.04%  '

$ws.Range("E35").Value = '  -2.65%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '59.38'
$ws.Range("E36").Value = '  -5.94%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  -0.07%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.395'
$ws.Range("E38").Value = '  -0.38%  '

$ws.Range("E39").Value = '  -4.61%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  +0.04%  '

$ws.Range("D41").Value = '0.0₃0713'
$ws.Range("E41").Value = '  +6.37%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.128'
$ws.Range("E42").Value = '  -1.99%  '

$ws.Range("D43").Value = '2.902.35'
$ws.Range("E43").Value = '  -4.94%  '

$ws.Range("E45").Value = '  -8.51%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0404'
$ws.Range("E46").Value = '  +2.31%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.69'
$ws.Range("E47").Value = '  -1.97%  '

$ws.Range("E48").Value = '  +12.35%  '

$ws.Range("E49").Value = '  -3.22%  '

$ws.Range("E50").Value = '  +0.78%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.99'
$ws.Range("E51").Value = '  +3.32%  '
